$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 8 with the new audit finding (mirrors the layout of rows 2-7)
$ws.Range("A8").Value = "Amélioration technique"
$ws.Range("B8").Value = "Mots clés en petit sur fond blanc"
$ws.Range("C8").Value = "Les div avec la classe ""keywords"" contient une succession de mots clés sans contexte. `nIl s'agit d'une technique de ""triche"" pour booster frauduleusement son référencement."
$ws.Range("D8").Value = "Essayer de tromper l'algorithme de Google représente un risque pour le référencement de la page car cela pourrait entrainer un malus.`nSuppression de ces divs là où elles apparaissent (header + footer)"
$ws.Range("E8").Value = "X"
$ws.Range("F8").Value = "smartkeyword - Black Hat"

# Match styling used by the other rows in the table
$ws.Range("A8").Style = $ws.Range("A6").Style
$ws.Range("B8").Style = $ws.Range("B7").Style
$ws.Range("C8").Style = $ws.Range("C7").Style
$ws.Range("D8").Style = $ws.Range("D7").Style
$ws.Range("E8").Style = $ws.Range("E7").Style
$ws.Range("F8").Style = $ws.Range("F7").Style

# Row height for the new content
$ws.Rows.Item(8).RowHeight = 52

# Column B is wider to fit the new content
$ws.Columns.Item(2).ColumnWidth = 34.140625

# Add the hyperlink for the new reference cell
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.smartkeyword.io/blog/black-hat-seo", "", "", "smartkeyword - Black Hat")

# Update the view: scroll back to A1, select F8
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F8").Select()
